# Update the "Förändrad" (Changed) date column (C) for all data rows (2-9)
# from 2026-02-21 (serial 46074) to 2026-02-22 (serial 46075).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 46075
}
